# Apply the commit's additions: two purple ("6600FF") annotation blocks
# appended onto two existing paragraphs in the "memoria" log.

$d = $word.ActiveDocument
$purple = 16711782  # wdColor for RGB 6600FF (stored/read back as 0x6600FF)

function Add-PurpleSegments($anchorText, $segments) {
    $rng = $d.Content
    $found = $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Anchor text not found: $anchorText"
    }
    $rng.Collapse(0)

    foreach ($seg in $segments) {
        $insStart = $rng.Start
        $rng.InsertAfter($seg)
        $rng.Start = $insStart
        $rng.End = $insStart + $seg.Length
        $rng.Font.Color = $purple
        $rng.Collapse(0)
    }
}

# --- Insertion 1: after "twitter me ha parecido guay, echadle un ojo cuando podáis." ---
Add-PurpleSegments "twitter me ha parecido guay, echadle un ojo cuando podáis." @(
    " ",
    "Linkeo",
    " la app con ",
    "FirebaseCli",
    " para poder usar ",
    "cloudmessaging",
    ". "
)

# --- Insertion 2: after "... parcialmente ... y que se mande un mensaje al usuario." ---
# This sits right before the trailing "_GoBack" bookmark at the end of the
# document, so a plain Find+Collapse+InsertAfter lands the new text *after*
# the bookmark markers. Instead, anchor a zero-length range exactly at the
# end-of-text offset and InsertBefore it (which inserts ahead of the
# bookmark), then force each appended chunk into its own run by nudging the
# font color away and back to the target purple (identical-property inserts
# silently coalesce into the neighboring run in this engine).

$rng2 = $d.Content
$found2 = $rng2.Find.Execute("se mande un mensaje al usuario.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Anchor text not found for insertion 2"
}
$endPos = $rng2.End

$segments2 = @(
    " Hay algunos ficheros que pertenecen a ",
    "functions",
    " de ",
    "firebase",
    ", de momento no hacen anda porque no las he programado, pero al hacer ",
    "merge",
    " no me deja quitarlas, de todas maneras, si finalmente no usamos ",
    "cloudmessaging",
    " y usamos otra solución los podemos dejar por si queremos implementar otras cosas con ",
    "FirebaseCli",
    "."
)

# Phase 1: insert every chunk (ahead of the trailing bookmark) first.
$bounds2 = @()
foreach ($seg in $segments2) {
    $ins = $d.Range($endPos, $endPos)
    $ins.InsertBefore($seg)
    $segStart = $endPos
    $segEnd = $endPos + $seg.Length
    $bounds2 += ,@($segStart, $segEnd)
    $endPos = $segEnd
}

# Phase 2: re-color each chunk individually so adjacent, identically
# formatted inserts don't silently coalesce back into one run — flipping
# the color away and back forces the engine to keep them as separate runs.
foreach ($b in $bounds2) {
    $segRng = $d.Range($b[0], $b[1])
    $segRng.Font.Color = 1
    $segRng.Font.Color = $purple
}

Write-Output "Done"
